# D_OPCO_DDL_mapping: repoint source from LND_CORE.D_DP_OPCO (table) to
# EDW.BI_D_OPCO_VW (view) as a ref source, and update the OPCO_ID /
# TBA_ACTIVE_FLG rows to pull straight from the view's columns instead of
# hard-coded constants. Also leaves the sheet selection on the Mapping tab
# (B5) as the last thing the author touched, matching the new DDL macro
# workflow.

$wb = $excel.ActiveWorkbook
$mapping = $wb.Worksheets.Item("Mapping")

# --- Header block: source table/type/name now point at the EDW view ---
$mapping.Range("B3").Value = "EBI_DEV_DB.EDW.BI_D_OPCO_VW"
$mapping.Range("B4").Value = "ref"
$mapping.Range("B5").Value = "EDW"

# --- Column mapping rows: every "Source Table" cell moves from the old
#     LND_CORE landing table to the new EDW view ---
$mapping.Range("C11").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C12").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C13").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C14").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C15").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C16").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C17").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C18").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C19").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C20").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C21").Value = "EDW.BI_D_OPCO_VW"

# OPCO_ID (row 10) and TBA_ACTIVE_FLG (row 22) used to be hard-coded
# constants ('1' / 'Y'); they now come straight from the view, so they gain
# a Source Table entry and their Logic/Mapping column switches from the
# quoted literal to the real column name. The leading "'" keeps the cells'
# existing text/quote-prefix formatting intact (same trick Excel uses when
# you retype a quote-prefixed cell). C10/C22 were blank before, so reset
# them back to the plain "Normal" style once they carry a value.
$mapping.Range("C10").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C10").Style = "Normal"
$mapping.Range("D10").Value = "'OPCO_ID"

$mapping.Range("C22").Value = "EDW.BI_D_OPCO_VW"
$mapping.Range("C22").Style = "Normal"
$mapping.Range("D22").Value = "'TBA_ACTIVE_FLG"

# --- Leave the workbook focused back on the Mapping sheet / cell B5 ---
$mapping.Activate()
$mapping.Range("B5").Select()
